$wb = $excel.ActiveWorkbook

# ----- Sheet: Summary -----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.400355871886121
$ws1.Range("C2").Value = 0.0743801652892562
$ws1.Range("D2").Value = 0.9642857142857143
$ws1.Range("E2").Value = 0.1381074168797954
$ws1.Range("F2").Value = 0.2842105263157895
$ws1.Range("G2").Value = 0.6603951081843838
$ws1.Range("H2").Value = 0.7998261102193687
$ws1.Range("I2").Value = 27
$ws1.Range("J2").Value = 336
$ws1.Range("K2").Value = 198
$ws1.Range("L2").Value = 1

# ----- Sheet: Classification Report -----
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 - class "0"
$ws2.Range("B2").Value = 0.9949748743718593
$ws2.Range("C2").Value = 0.3707865168539326
$ws2.Range("D2").Value = 0.5402455661664393

# Row 3 - class "1"
$ws2.Range("B3").Value = 0.0743801652892562
$ws2.Range("C3").Value = 0.9642857142857143
$ws2.Range("D3").Value = 0.1381074168797954

# Row 4 - accuracy
$ws2.Range("B4").Value = 0.400355871886121
$ws2.Range("C4").Value = 0.400355871886121
$ws2.Range("D4").Value = 0.400355871886121
$ws2.Range("E4").Value = 0.400355871886121

# Row 5 - macro avg
$ws2.Range("B5").Value = 0.5346775198305578
$ws2.Range("C5").Value = 0.6675361155698234
$ws2.Range("D5").Value = 0.3391764915231174

# Row 6 - weighted avg
$ws2.Range("B6").Value = 0.9491089458054663
$ws2.Range("C6").Value = 0.400355871886121
$ws2.Range("D6").Value = 0.5202102135329409

# ----- Sheet: Confusion Matrix -----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - Actual 0
$ws3.Range("B2").Value = 198
$ws3.Range("C2").Value = 336

# Row 3 - Actual 1
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 27
